$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new blank rows before row 48 (shifts old rows 48-65 down to 54-71)
$ws.Rows.Item(48).Resize(6).Insert()

# Remove the leftover column-D formatting on rows that should stay fully empty
$ws.Range("D49").Clear()
$ws.Range("D51").Clear()
$ws.Range("D52").Clear()
$ws.Range("D53").Clear()

# New data row: date, hours and task note
$ws.Range("A48").Value = 43548
$ws.Range("B48").Value = 2.5
$ws.Range("D48").Value = "Project/Presentation: Got code working to backfill the indexes for Hibernate Search; completed a successful search; worked on PowerPoint"
$ws.Rows.Item(48).RowHeight = 30

# New note (leading apostrophe forces text / quotePrefix since it starts with a number)
$ws.Range("D50").Value = "'2.5 doucmented plus Sun 10:10 - 11:45"

# Restore the view the author left the sheet in
$ws.Application.ActiveWindow.ScrollRow = 46
